$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -18.42350578835206
$ws.Cells.Item(2, 3).Value = 1.937749438846114
$ws.Cells.Item(2, 4).Value = -18.42350578835206
$ws.Cells.Item(2, 5).Value = -18.42350578835206
$ws.Cells.Item(2, 6).Value = -18.42350578835206
$ws.Cells.Item(2, 7).Value = -18.42350578835206
$ws.Cells.Item(2, 8).Value = -18.42350578835206
$ws.Cells.Item(2, 9).Value = -18.42350578835206
$ws.Cells.Item(2, 10).Value = -18.42350578835206
$ws.Cells.Item(2, 11).Value = -18.42350578835206
$ws.Cells.Item(3, 2).Value = -18.42350578835206
$ws.Cells.Item(3, 3).Value = -18.42350578835206
$ws.Cells.Item(3, 4).Value = -18.42350578835206
$ws.Cells.Item(3, 5).Value = -18.42350578835206
$ws.Cells.Item(3, 6).Value = -18.42350578835206
$ws.Cells.Item(3, 7).Value = -18.42350578835206
$ws.Cells.Item(3, 8).Value = -18.42350578835206
$ws.Cells.Item(3, 9).Value = 4.321924196626807
$ws.Cells.Item(3, 10).Value = -18.42350578835206
$ws.Cells.Item(3, 11).Value = -18.42350578835206
$ws.Cells.Item(4, 2).Value = -18.42350578835206
$ws.Cells.Item(4, 3).Value = 1.988556440912792
$ws.Cells.Item(4, 4).Value = 1.541124160370487
$ws.Cells.Item(4, 5).Value = -18.42350578835206
$ws.Cells.Item(4, 6).Value = 3.378700131825844
$ws.Cells.Item(4, 7).Value = -18.42350578835206
$ws.Cells.Item(4, 8).Value = 0.9358856099001401
$ws.Cells.Item(4, 9).Value = -18.42350578835206
$ws.Cells.Item(4, 10).Value = -0.7738776999247209
$ws.Cells.Item(4, 11).Value = -18.42350578835206
$ws.Cells.Item(5, 2).Value = -18.42350578835206
$ws.Cells.Item(5, 3).Value = 1.72454992984049
$ws.Cells.Item(5, 4).Value = -18.42350578835206
$ws.Cells.Item(5, 5).Value = -18.42350578835206
$ws.Cells.Item(5, 6).Value = -18.42350578835206
$ws.Cells.Item(5, 7).Value = 2.570392452606983
$ws.Cells.Item(5, 8).Value = -18.42350578835206
$ws.Cells.Item(5, 9).Value = -18.42350578835206
$ws.Cells.Item(5, 10).Value = -18.42350578835206
$ws.Cells.Item(5, 11).Value = -18.42350578835206
$ws.Cells.Item(6, 2).Value = -18.42350578835206
$ws.Cells.Item(6, 3).Value = -18.42350578835206
$ws.Cells.Item(6, 4).Value = -18.42350578835206
$ws.Cells.Item(6, 5).Value = -18.42350578835206
$ws.Cells.Item(6, 6).Value = -18.42350578835206
$ws.Cells.Item(6, 7).Value = -18.42350578835206
$ws.Cells.Item(6, 8).Value = -18.42350578835206
$ws.Cells.Item(6, 9).Value = -18.42350578835206
$ws.Cells.Item(6, 10).Value = -18.42350578835206
$ws.Cells.Item(6, 11).Value = -18.42350578835206
$ws.Cells.Item(7, 2).Value = 2.437521276795743
$ws.Cells.Item(7, 3).Value = -18.42350578835206
$ws.Cells.Item(7, 4).Value = -18.42350578835206
$ws.Cells.Item(7, 5).Value = -18.42350578835206
$ws.Cells.Item(7, 6).Value = -18.42350578835206
$ws.Cells.Item(7, 7).Value = -18.42350578835206
$ws.Cells.Item(7, 8).Value = -18.42350578835206
$ws.Cells.Item(7, 9).Value = -18.42350578835206
$ws.Cells.Item(7, 10).Value = -18.42350578835206
$ws.Cells.Item(7, 11).Value = -18.42350578835206
$ws.Cells.Item(8, 2).Value = -18.42350578835206
$ws.Cells.Item(8, 3).Value = -18.42350578835206
$ws.Cells.Item(8, 4).Value = -18.42350578835206
$ws.Cells.Item(8, 5).Value = 1.784286663620522
$ws.Cells.Item(8, 6).Value = -18.42350578835206
$ws.Cells.Item(8, 7).Value = -18.42350578835206
$ws.Cells.Item(8, 8).Value = -18.42350578835206
$ws.Cells.Item(8, 9).Value = -18.42350578835206
$ws.Cells.Item(8, 10).Value = -18.42350578835206
$ws.Cells.Item(8, 11).Value = -18.42350578835206
$ws.Cells.Item(9, 2).Value = 3.86620047125174
$ws.Cells.Item(9, 3).Value = -18.42350578835206
$ws.Cells.Item(9, 4).Value = -18.42350578835206
$ws.Cells.Item(9, 5).Value = -18.42350578835206
$ws.Cells.Item(9, 6).Value = -18.42350578835206
$ws.Cells.Item(9, 7).Value = -18.42350578835206
$ws.Cells.Item(9, 8).Value = -18.42350578835206
$ws.Cells.Item(9, 9).Value = -18.42350578835206
$ws.Cells.Item(9, 10).Value = -18.42350578835206
$ws.Cells.Item(9, 11).Value = -18.42350578835206
$ws.Cells.Item(10, 2).Value = -18.42350578835206
$ws.Cells.Item(10, 3).Value = -18.42350578835206
$ws.Cells.Item(10, 4).Value = -18.42350578835206
$ws.Cells.Item(10, 5).Value = -18.42350578835206
$ws.Cells.Item(10, 6).Value = -18.42350578835206
$ws.Cells.Item(10, 7).Value = -18.42350578835206
$ws.Cells.Item(10, 8).Value = -18.42350578835206
$ws.Cells.Item(10, 9).Value = -18.42350578835206
$ws.Cells.Item(10, 10).Value = -18.42350578835206
$ws.Cells.Item(10, 11).Value = 2.242927147341175
$ws.Cells.Item(11, 2).Value = -18.42350578835206
$ws.Cells.Item(11, 3).Value = -18.42350578835206
$ws.Cells.Item(11, 4).Value = -18.42350578835206
$ws.Cells.Item(11, 5).Value = 2.915416076012164
$ws.Cells.Item(11, 6).Value = -18.42350578835206
$ws.Cells.Item(11, 7).Value = 3.032016595051148
$ws.Cells.Item(11, 8).Value = -18.42350578835206
$ws.Cells.Item(11, 9).Value = -18.42350578835206
$ws.Cells.Item(11, 10).Value = -18.42350578835206
$ws.Cells.Item(11, 11).Value = 1.84437016594062
$ws.Cells.Item(12, 2).Value = -18.42350578835206
$ws.Cells.Item(12, 3).Value = -18.42350578835206
$ws.Cells.Item(12, 4).Value = -18.42350578835206
$ws.Cells.Item(12, 5).Value = -18.42350578835206
$ws.Cells.Item(12, 6).Value = -18.42350578835206
$ws.Cells.Item(12, 7).Value = -18.42350578835206
$ws.Cells.Item(12, 8).Value = -18.42350578835206
$ws.Cells.Item(12, 9).Value = -18.42350578835206
$ws.Cells.Item(12, 10).Value = -18.42350578835206
$ws.Cells.Item(12, 11).Value = -18.42350578835206
$ws.Cells.Item(13, 2).Value = -18.42350578835206
$ws.Cells.Item(13, 3).Value = -18.42350578835206
$ws.Cells.Item(13, 4).Value = -18.42350578835206
$ws.Cells.Item(13, 5).Value = 2.554805967722549
$ws.Cells.Item(13, 6).Value = -18.42350578835206
$ws.Cells.Item(13, 7).Value = -18.42350578835206
$ws.Cells.Item(13, 8).Value = -18.42350578835206
$ws.Cells.Item(13, 9).Value = -18.42350578835206
$ws.Cells.Item(13, 10).Value = 1.323443104926883
$ws.Cells.Item(13, 11).Value = 1.635116225086647
$ws.Cells.Item(14, 2).Value = -18.42350578835206
$ws.Cells.Item(14, 3).Value = -18.42350578835206
$ws.Cells.Item(14, 4).Value = 1.535869366348396
$ws.Cells.Item(14, 5).Value = -18.42350578835206
$ws.Cells.Item(14, 6).Value = -18.42350578835206
$ws.Cells.Item(14, 7).Value = -18.42350578835206
$ws.Cells.Item(14, 8).Value = -18.42350578835206
$ws.Cells.Item(14, 9).Value = -18.42350578835206
$ws.Cells.Item(14, 10).Value = -18.42350578835206
$ws.Cells.Item(14, 11).Value = 2.038094557915746
$ws.Cells.Item(15, 2).Value = -18.42350578835206
$ws.Cells.Item(15, 3).Value = -18.42350578835206
$ws.Cells.Item(15, 4).Value = 1.875099415755087
$ws.Cells.Item(15, 5).Value = -18.42350578835206
$ws.Cells.Item(15, 6).Value = -18.42350578835206
$ws.Cells.Item(15, 7).Value = -18.42350578835206
$ws.Cells.Item(15, 8).Value = -18.42350578835206
$ws.Cells.Item(15, 9).Value = -18.42350578835206
$ws.Cells.Item(15, 10).Value = -18.42350578835206
$ws.Cells.Item(15, 11).Value = -18.42350578835206
$ws.Cells.Item(16, 2).Value = -18.42350578835206
$ws.Cells.Item(16, 3).Value = -18.42350578835206
$ws.Cells.Item(16, 4).Value = -18.42350578835206
$ws.Cells.Item(16, 5).Value = -18.42350578835206
$ws.Cells.Item(16, 6).Value = -18.42350578835206
$ws.Cells.Item(16, 7).Value = -18.42350578835206
$ws.Cells.Item(16, 8).Value = -18.42350578835206
$ws.Cells.Item(16, 9).Value = -18.42350578835206
$ws.Cells.Item(16, 10).Value = 1.950099647865605
$ws.Cells.Item(16, 11).Value = -18.42350578835206
$ws.Cells.Item(17, 2).Value = -18.42350578835206
$ws.Cells.Item(17, 3).Value = 2.107952777828027
$ws.Cells.Item(17, 4).Value = 1.885420896937866
$ws.Cells.Item(17, 5).Value = -18.42350578835206
$ws.Cells.Item(17, 6).Value = -18.42350578835206
$ws.Cells.Item(17, 7).Value = -18.42350578835206
$ws.Cells.Item(17, 8).Value = 2.031545449349045
$ws.Cells.Item(17, 9).Value = -18.42350578835206
$ws.Cells.Item(17, 10).Value = 2.845740195553066
$ws.Cells.Item(17, 11).Value = -18.42350578835206
$ws.Cells.Item(18, 2).Value = -18.42350578835206
$ws.Cells.Item(18, 3).Value = -18.42350578835206
$ws.Cells.Item(18, 4).Value = -18.42350578835206
$ws.Cells.Item(18, 5).Value = -18.42350578835206
$ws.Cells.Item(18, 6).Value = -18.42350578835206
$ws.Cells.Item(18, 7).Value = -18.42350578835206
$ws.Cells.Item(18, 8).Value = 2.305773895778891
$ws.Cells.Item(18, 9).Value = -18.42350578835206
$ws.Cells.Item(18, 10).Value = 2.550833873881191
$ws.Cells.Item(18, 11).Value = -18.42350578835206
$ws.Cells.Item(19, 2).Value = -18.42350578835206
$ws.Cells.Item(19, 3).Value = -18.42350578835206
$ws.Cells.Item(19, 4).Value = 1.992232885757828
$ws.Cells.Item(19, 5).Value = -18.42350578835206
$ws.Cells.Item(19, 6).Value = -18.42350578835206
$ws.Cells.Item(19, 7).Value = -18.42350578835206
$ws.Cells.Item(19, 8).Value = 1.562259727985106
$ws.Cells.Item(19, 9).Value = -18.42350578835206
$ws.Cells.Item(19, 10).Value = -18.42350578835206
$ws.Cells.Item(19, 11).Value = -18.42350578835206
$ws.Cells.Item(20, 2).Value = -18.42350578835206
$ws.Cells.Item(20, 3).Value = 1.065960018604508
$ws.Cells.Item(20, 4).Value = 1.510313072035357
$ws.Cells.Item(20, 5).Value = -18.42350578835206
$ws.Cells.Item(20, 6).Value = 3.262822465379656
$ws.Cells.Item(20, 7).Value = -18.42350578835206
$ws.Cells.Item(20, 8).Value = 1.758540546230999
$ws.Cells.Item(20, 9).Value = -18.42350578835206
$ws.Cells.Item(20, 10).Value = -18.42350578835206
$ws.Cells.Item(20, 11).Value = 2.157780317157008
$ws.Cells.Item(21, 2).Value = -18.42350578835206
$ws.Cells.Item(21, 3).Value = 1.316935443984572
$ws.Cells.Item(21, 4).Value = -18.42350578835206
$ws.Cells.Item(21, 5).Value = 1.648528655354912
$ws.Cells.Item(21, 6).Value = -18.42350578835206
$ws.Cells.Item(21, 7).Value = 2.555992675350175
$ws.Cells.Item(21, 8).Value = 1.442241072936099
$ws.Cells.Item(21, 9).Value = -18.42350578835206
$ws.Cells.Item(21, 10).Value = -18.42350578835206
$ws.Cells.Item(21, 11).Value = -18.42350578835206
